# Update column G (header "K", strikeouts) values for rows 2-41
# per commit: "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(2,7,6,4,1,6,6,8,9,5,11,9,4,3,10,11,4,1,9,6,5,3,6,7,3,10,6,4,6,1,3,6,3,4,6,8,7,6,5,3)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
